$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-25 Friday", "2025-07-26 Saturday"),
    @("92×47=4324", "52×93=4836"),
    @("23×65=1495", "78×11=858"),
    @("52×25=1300", "52×31=1612"),
    @("71×77=5467", "42×82=3444"),
    @("67×63=4221", "46×49=2254"),
    @("88×18=1584", "86×63=5418"),
    @("41×57=2337", "93×49=4557"),
    @("88×44=3872", "39×62=2418"),
    @("33×15=495", "52×72=3744"),
    @("38×94=3572", "57×98=5586"),
    @("55×44=2420", "19×86=1634"),
    @("72×87=6264", "55×84=4620"),
    @("67×51=3417", "66×38=2508"),
    @("38×77=2926", "95×81=7695"),
    @("55×71=3905", "48×89=4272"),
    @("75×45=3375", "28×15=420"),
    @("44×44=1936", "46×26=1196"),
    @("57×14=798", "74×12=888"),
    @("19×11=209", "65×12=780"),
    @("92×17=1564", "31×22=682"),
    @("64×36=2304", "52×18=936"),
    @("65×46=2990", "73×54=3942"),
    @("35×32=1120", "63×98=6174"),
    @("75×86=6450", "87×54=4698"),
    @("28×37=1036", "71×79=5609")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
